$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(164.1, -64.59999999999999, -9.119999999999999, 979, 3),
    @(139.4, -56.4, -2.67, 700.78, 3),
    @(138, -55.8, -2.15, 700.67, 3),
    @(164, -64, -9.17, 979.89, 3),
    @(80.43000000000001, -50, 66.89, 521, 3),
    @(109.723, -54, 50.661, 610, 3),
    @(138, -55.8, -2.15, 700.67, 3),
    @(161.443, -58.8, -5.3, 839.78, 3),
    @(159.7, -59.2, -5.48, 839, 3),
    @(109, -55.8, 50.661, 610.67, 3)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
